# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from CompressionAlgor" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from CompressionAlgor")
$wsInclude.Name = "Include #0"

# --- 2. Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
# Copy formatting for the newly-extended last row (15) from the current last row (14)
# first, then shift the row 11..14 contents down to 12..15, and finally
# write the new Jurisdiction row into row 11. This keeps the existing
# style indices intact rather than generating new ones (as Rows.Insert would).
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

for ($r = 14; $r -ge 11; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Text
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Text
}

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
